$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Match the author's resized/repositioned workbook window ---
$win = $excel.ActiveWindow
$win.Left = 15360
$win.Top = 780
$win.Width = 21285
$win.Height = 13215

# --- Update ingredient / allergen / diet / label text for the burger, mac & cheese,
#     and mozza sticks rows (B2, B3, B4, C4, E4, F4, B5, F5) ---

# Cheese Burger (row 2): strip the bracketed beef-patty sourcing detail
$ws.Range("B2").Value = "Sesame Burger Bun / Beef Patty / Pickles / Cheddar Cheese / Red Onion"

# Veggie Burger (row 3): strip the bracketed black-bean patty sourcing detail
$ws.Range("B3").Value = "Sesame Burger Bun / Spicy Black Bean Patty / Pickles / Cheddar Cheese / Red Onion"

# Mac & Cheese Bites (row 4): fill in real allergens / ingredients / diet / nutrition label
$ws.Range("C4").Value = "Wheat, gluten, milk, eggs."
$ws.Range("B4").Value = "Macaroni / Cream / Flour / Cheese Blend / Egg / Breading"
$ws.Range("E4").Value = "VEG"

# Mozza Sticks (row 5): fill in real ingredients
$ws.Range("B5").Value = "Mozzarella / Flour / Eggs / Seasoned Breading / Salt"

# Nutrition label codes for the two new rows
$ws.Range("F4").Value = "Mac_&_Cheese_Bites"
$ws.Range("F5").Value = "Mozza_Sticks"

# --- Update the selected cell shown when the sheet is opened ---
$null = $ws.Range("F5").Select()
